# Mark additional "h0=uncond" / "h0=ht MLE P" scenario runs as complete
# (status value 1) in the "ongoing" mse columns of the lower status table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# "h0=uncond / ongoing" block (column K, mse) - years 2011-2015
$ws.Range("K23").Value = 1   # 2011
$ws.Range("K24").Value = 1   # 2012
$ws.Range("K25").Value = 1   # 2013
$ws.Range("K26").Value = 1   # 2014
$ws.Range("K27").Value = 1   # 2015

# "h0=ht MLE P / ongoing" block (column C, mse) - years 2014-2018
$ws.Range("C26").Value = 1   # 2014
$ws.Range("C27").Value = 1   # 2015
$ws.Range("C28").Value = 1   # 2016
$ws.Range("C29").Value = 1   # 2017
$ws.Range("C30").Value = 1   # 2018

# Leave the cursor where the author last left it when saving.
$ws.Range("K29").Select()
